# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Leve profit sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) per sheet.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 2100.1667
$ws.Range("I98").Value = 2020.2
$ws.Range("K98").Value = 2020.2
$ws.Range("M98").Value = -522.2
# Row 122
$ws.Range("H122").Value = 2100.1667
$ws.Range("I122").Value = 2020.2
$ws.Range("K122").Value = 6060.6
$ws.Range("M122").Value = -3610.6
# Row 132
$ws.Range("H132").Value = 2668.5
$ws.Range("I132").Value = 2518.4048
$ws.Range("J132").Value = 4244.5
$ws.Range("K132").Value = 7555.214399999999
$ws.Range("L132").Value = 12733.5
$ws.Range("M132").Value = -5025.214399999999
$ws.Range("N132").Value = -17793.5
# Row 135
$ws.Range("H135").Value = 50001804
$ws.Range("I135").Value = 83335176
$ws.Range("J135").Value = 1743
$ws.Range("K135").Value = 750016584
$ws.Range("L135").Value = 15687
$ws.Range("M135").Value = -750014049
$ws.Range("N135").Value = -20757
# Row 138
$ws.Range("H138").Value = 4192.0273
$ws.Range("I138").Value = 3763.5
$ws.Range("J138").Value = 4332.273
$ws.Range("K138").Value = 11290.5
$ws.Range("L138").Value = 12996.819
$ws.Range("M138").Value = -6150.5
$ws.Range("N138").Value = -23276.819

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 11114705
$ws.Range("I61").Value = 27780472
$ws.Range("J61").Value = 4193.1665
$ws.Range("K61").Value = 27780472
$ws.Range("L61").Value = 4193.1665
$ws.Range("M61").Value = -27780260
$ws.Range("N61").Value = -4617.1665
# Row 62
$ws.Range("H62").Value = 23666.666
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 33000
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 33000
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -34248
# Row 65
$ws.Range("H65").Value = 23666.666
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 33000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 99000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -105240
# Row 102
$ws.Range("H102").Value = 3512.8572
$ws.Range("I102").Value = 3615
$ws.Range("K102").Value = 3615
$ws.Range("M102").Value = -1993
# Row 108
$ws.Range("H108").Value = 75000
$ws.Range("J108").Value = 75000
$ws.Range("L108").Value = 75000
$ws.Range("N108").Value = -82680
# Row 115
$ws.Range("H115").Value = 62500
$ws.Range("J115").Value = 62500
$ws.Range("L115").Value = 62500
$ws.Range("N115").Value = -65634
# Row 136
$ws.Range("H136").Value = 11114705
$ws.Range("I136").Value = 27780472
$ws.Range("J136").Value = 4193.1665
$ws.Range("K136").Value = 83341416
$ws.Range("L136").Value = 12579.4995
$ws.Range("M136").Value = -83338866
$ws.Range("N136").Value = -17679.4995

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 103
$ws.Range("H103").Value = 98657
$ws.Range("J103").Value = 98657
$ws.Range("L103").Value = 98657
$ws.Range("N103").Value = -101001
# Row 109
$ws.Range("H109").Value = 67342
$ws.Range("J109").Value = 67342
$ws.Range("L109").Value = 67342
$ws.Range("N109").Value = -70116
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# Row 134
$ws.Range("H134").Value = 3368.1956
$ws.Range("I134").Value = 3331.718
$ws.Range("J134").Value = 3571.4285
$ws.Range("K134").Value = 9995.153999999999
$ws.Range("L134").Value = 10714.2855
$ws.Range("M134").Value = -7460.153999999999
$ws.Range("N134").Value = -15784.2855

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8228.365
$ws.Range("I31").Value = 3196.3
$ws.Range("J31").Value = 9177.812
$ws.Range("K31").Value = 3196.3
$ws.Range("L31").Value = 9177.812
$ws.Range("M31").Value = -2901.3
$ws.Range("N31").Value = -9767.812
# Row 34
$ws.Range("H34").Value = 8228.365
$ws.Range("I34").Value = 3196.3
$ws.Range("J34").Value = 9177.812
$ws.Range("K34").Value = 3196.3
$ws.Range("L34").Value = 9177.812
$ws.Range("M34").Value = -2994.3
$ws.Range("N34").Value = -9581.812
# Row 114
$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
# Row 134
$ws.Range("H134").Value = 4633079
$ws.Range("I134").Value = 5438387.5
$ws.Range("J134").Value = 2555.25
$ws.Range("K134").Value = 16315162.5
$ws.Range("L134").Value = 7665.75
$ws.Range("M134").Value = -16312627.5
$ws.Range("N134").Value = -12735.75

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 708.3889
$ws.Range("I113").Value = 699.2105
$ws.Range("K113").Value = 2097.6315
$ws.Range("M113").Value = 72.36850000000004
# Row 127
$ws.Range("H127").Value = 452.2857
$ws.Range("J127").Value = 452.2857
$ws.Range("L127").Value = 1356.8571
$ws.Range("N127").Value = -11276.8571
# Row 132
$ws.Range("H132").Value = 2656.6047
$ws.Range("J132").Value = 3043.8667
$ws.Range("L132").Value = 27394.8003
$ws.Range("N132").Value = -32454.8003

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
# Row 99
$ws.Range("H99").Value = 9718.611000000001
$ws.Range("I99").Value = 3841.1538
$ws.Range("J99").Value = 25000
$ws.Range("K99").Value = 3841.1538
$ws.Range("L99").Value = 25000
$ws.Range("M99").Value = -1595.1538
$ws.Range("N99").Value = -29492
# Row 107
$ws.Range("H107").Value = 465.08334
$ws.Range("I107").Value = 398.2
$ws.Range("J107").Value = 799.5
$ws.Range("K107").Value = 398.2
$ws.Range("L107").Value = 799.5
$ws.Range("M107").Value = 1521.8
$ws.Range("N107").Value = -4639.5
# Row 132
$ws.Range("H132").Value = 26320944
$ws.Range("I132").Value = 45461388
$ws.Range("J132").Value = 2832.9375
$ws.Range("K132").Value = 136384164
$ws.Range("L132").Value = 8498.8125
$ws.Range("M132").Value = -136381634
$ws.Range("N132").Value = -13558.8125
# Row 133
$ws.Range("H133").Value = 52740
$ws.Range("J133").Value = 52740
$ws.Range("L133").Value = 52740
$ws.Range("N133").Value = -62860

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
# Row 117
$ws.Range("H117").Value = 98392
$ws.Range("J117").Value = 98392
$ws.Range("L117").Value = 98392
$ws.Range("N117").Value = -107570
# Row 132
$ws.Range("H132").Value = 3403.0278
$ws.Range("I132").Value = 2866.3333
$ws.Range("J132").Value = 3786.3809
$ws.Range("K132").Value = 8598.999899999999
$ws.Range("L132").Value = 11359.1427
$ws.Range("M132").Value = -6068.999899999999
$ws.Range("N132").Value = -16419.1427
# Row 135
$ws.Range("H135").Value = 89000
$ws.Range("J135").Value = 89000
$ws.Range("L135").Value = 89000
$ws.Range("N135").Value = -99140
# Row 140
$ws.Range("H140").Value = 73399.60000000001
$ws.Range("J140").Value = 73399.60000000001
$ws.Range("L140").Value = 73399.60000000001
$ws.Range("N140").Value = -83759.60000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 386674
$ws.Range("J29").Value = 386674
$ws.Range("L29").Value = 386674
$ws.Range("N29").Value = -387254
# Row 46
$ws.Range("H46").Value = 240999
$ws.Range("J46").Value = 240999
$ws.Range("L46").Value = 240999
$ws.Range("N46").Value = -241461
# Row 76
$ws.Range("H76").Value = 41868.6
$ws.Range("J76").Value = 41868.6
$ws.Range("L76").Value = 41868.6
$ws.Range("N76").Value = -42498.6
# Row 79
$ws.Range("H79").Value = 41868.6
$ws.Range("J79").Value = 41868.6
$ws.Range("L79").Value = 41868.6
$ws.Range("N79").Value = -44052.6
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 129
$ws.Range("H129").Value = 98429
$ws.Range("J129").Value = 98429
$ws.Range("L129").Value = 98429
$ws.Range("N129").Value = -108429
# Row 134
$ws.Range("H134").Value = 240999
$ws.Range("J134").Value = 240999
$ws.Range("L134").Value = 722997
$ws.Range("N134").Value = -728067
